# Apply the "杭州-漫展信息" update (generated output at commit 456a3b4).
#
# Sheet map (per xl/workbook.xml order, 1-based):
#   1 = 展览    (Exhibitions)
#   2 = 演出    (Performances)
#   3 = 本地生活 (Local life)
#   4 = 全部类型 (All types / combined feed)
#
# Helper: write a literal text value into a cell while preventing Excel's
# COM layer from auto-coercing date-looking strings ("2024-05-12") into
# date serial numbers -- the source file stores these as plain text.
function Set-TextCell($ws, [string]$addr, [string]$text) {
    $ws.Range($addr).Value = "'" + $text
}

function Set-NumCell($ws, [string]$addr, $num) {
    $ws.Range($addr).Value = $num
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (index 1): simple "想去人数" (F column) count bumps.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

Set-NumCell $ws1 "F3"  774
Set-NumCell $ws1 "F4"  508
Set-NumCell $ws1 "F7"  1097
Set-NumCell $ws1 "F9"  19
Set-NumCell $ws1 "F15" 726
Set-NumCell $ws1 "F16" 782
Set-NumCell $ws1 "F19" 49
Set-NumCell $ws1 "F20" 619
Set-NumCell $ws1 "F21" 108
Set-NumCell $ws1 "F23" 1911
Set-NumCell $ws1 "F24" 484
Set-NumCell $ws1 "F26" 1726
Set-NumCell $ws1 "F27" 257
Set-NumCell $ws1 "F28" 2536
Set-NumCell $ws1 "F29" 453
Set-NumCell $ws1 "F31" 651
Set-NumCell $ws1 "F35" 875
Set-NumCell $ws1 "F36" 1584
Set-NumCell $ws1 "F37" 266
Set-NumCell $ws1 "F40" 109
Set-NumCell $ws1 "F41" 97
Set-NumCell $ws1 "F42" 136

# ---------------------------------------------------------------------
# Sheet "全部类型" (index 4)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# Same simple F-count bumps as above, rows are offset by +2 on this sheet.
Set-NumCell $ws4 "F5"  774
Set-NumCell $ws4 "F6"  508
Set-NumCell $ws4 "F9"  1097
Set-NumCell $ws4 "F11" 19
Set-NumCell $ws4 "F16" 726
Set-NumCell $ws4 "F17" 782
Set-NumCell $ws4 "F24" 49
Set-NumCell $ws4 "F25" 619
Set-NumCell $ws4 "F26" 108
Set-NumCell $ws4 "F28" 1911
Set-NumCell $ws4 "F29" 484
Set-NumCell $ws4 "F32" 2536
Set-NumCell $ws4 "F33" 453

# A new event ("杭州·《卡农》永恒经典名曲音乐会") was published for 2024-05-12,
# sorted ahead of the existing "奇迹の闪耀" row. This pushes rows 35-43 down
# by one; the previously-last of them ("苗阜王声 ...") drops off, and row 45
# onward is unaffected (it absorbs the shift).

# New row 35: 杭州·《卡农》永恒经典名曲音乐会 (B35 already reads 2024-05-12)
Set-TextCell $ws4 "C35" "杭州·《卡农》永恒经典名曲音乐会"
Set-TextCell $ws4 "D35" "武林路77号 浙江省文化馆小剧场（原群艺馆小剧场）"
Set-TextCell $ws4 "E35" "2024.05.12 14:00-05.12 15:30"
Set-NumCell  $ws4 "F35" 3
Set-NumCell  $ws4 "G35" 100
Set-TextCell $ws4 "H35" "https://show.bilibili.com/platform/detail.html?id=83176"
Set-TextCell $ws4 "I35" "//i0.hdslb.com/bfs/openplatform/202403/gLrSkh0O1711013683966.jpeg"

# Row 36 <- old row 35 (杭州·奇迹の闪耀 「UP!」巡回动漫演唱会)
Set-TextCell $ws4 "B36" "2024-05-12"
Set-TextCell $ws4 "C36" "杭州·奇迹の闪耀 「UP!」巡回动漫演唱会"
Set-TextCell $ws4 "D36" "东坡路10号 杭州东坡大剧院"
Set-TextCell $ws4 "E36" "2024.05.12 19:30-05.12 21:30"
Set-NumCell  $ws4 "F36" 12
Set-NumCell  $ws4 "G36" 126
Set-TextCell $ws4 "H36" "https://show.bilibili.com/platform/detail.html?id=82452"
Set-TextCell $ws4 "I36" "//i1.hdslb.com/bfs/openplatform/202403/HvxHPz981709707512970.jpeg"

# Row 37 <- old row 36 (杭州·Jo迪)
Set-TextCell $ws4 "C37" "杭州·Jo迪"
Set-TextCell $ws4 "D37" "萧杭路28号 格拉斯club"
Set-TextCell $ws4 "E37" "2024.05.18 13:00-05.18 19:00"
Set-NumCell  $ws4 "F37" 33
Set-NumCell  $ws4 "G37" 198
Set-TextCell $ws4 "H37" "https://show.bilibili.com/platform/detail.html?id=83008"
Set-TextCell $ws4 "I37" "//i1.hdslb.com/bfs/openplatform/202403/AEtl5BHN1711015003341.jpeg"

# Row 38 <- old row 37 (杭州·《沐云华·次元狂想》经典动漫二次元ACG音乐会)
Set-TextCell $ws4 "C38" "杭州·《沐云华·次元狂想》经典动漫二次元ACG音乐会"
Set-TextCell $ws4 "D38" "建国南路280号 杭州红星剧院"
Set-TextCell $ws4 "E38" "2024.05.18 19:30-05.18 22:00"
Set-NumCell  $ws4 "F38" 60
Set-NumCell  $ws4 "G38" 90
Set-TextCell $ws4 "H38" "https://show.bilibili.com/platform/detail.html?id=83113"
Set-TextCell $ws4 "I38" "//i1.hdslb.com/bfs/openplatform/202403/TXmgAvCC1710582339525.jpeg"

# Row 39 <- old row 38 (杭州·现世繁华-代号鸢only), F bumped 650 -> 651
Set-TextCell $ws4 "B39" "2024-05-18"
Set-TextCell $ws4 "C39" "杭州·现世繁华-代号鸢only"
Set-TextCell $ws4 "D39" "石祥路575号 杭州海外海纳川大酒店(万达广场渡驾桥地铁站店)"
Set-TextCell $ws4 "E39" "2024.05.18 10:00-05.18 21:00"
Set-NumCell  $ws4 "F39" 651
Set-NumCell  $ws4 "G39" 76
Set-TextCell $ws4 "H39" "https://show.bilibili.com/platform/detail.html?id=81905"
Set-TextCell $ws4 "I39" "//i2.hdslb.com/bfs/openplatform/202402/m3upuV2F1708327958926.jpeg"

# Row 40 <- old row 39 (杭州·原神X星铁X绝区零only)
Set-TextCell $ws4 "C40" "杭州·原神X星铁X绝区零only"
Set-TextCell $ws4 "D40" "望江东路333号 杭州瑞莱克斯大酒店"
Set-TextCell $ws4 "E40" "2024.05.25 10:00-05.25 17:00"
Set-NumCell  $ws4 "F40" 125
Set-NumCell  $ws4 "G40" 60
Set-TextCell $ws4 "H40" "https://show.bilibili.com/platform/detail.html?id=82754"
Set-TextCell $ws4 "I40" "//i1.hdslb.com/bfs/openplatform/202403/qA0LNJuF1710234461030.jpeg"

# Row 41 <- old row 40 (杭州·第三届缘起cp展 我们二次元的情人节！)
Set-TextCell $ws4 "B41" "2024-05-25"
Set-TextCell $ws4 "C41" "杭州·第三届缘起cp展 我们二次元的情人节！"
Set-TextCell $ws4 "D41" "黄姑山路51-4号 0101park"
Set-TextCell $ws4 "E41" "2024.05.25 10:00-05.26 17:00"
Set-NumCell  $ws4 "F41" 84
Set-NumCell  $ws4 "G41" 65
Set-TextCell $ws4 "H41" "https://show.bilibili.com/platform/detail.html?id=83336"
Set-TextCell $ws4 "I41" "//i1.hdslb.com/bfs/openplatform/202403/D9t8ms7G1711350634757.png"

# Row 42 <- old row 41 (杭州·造梦探险家——二次元同好会)
Set-TextCell $ws4 "B42" "2024-06-01"
Set-TextCell $ws4 "C42" "杭州·造梦探险家——二次元同好会"
Set-TextCell $ws4 "D42" "临平街道北沙西路156-1号 杭州临平遇上设计师酒店"
Set-TextCell $ws4 "E42" "2024.06.01 10:00-06.01 16:00"
Set-NumCell  $ws4 "F42" 88
Set-NumCell  $ws4 "G42" 28
Set-TextCell $ws4 "H42" "https://show.bilibili.com/platform/detail.html?id=82736"
Set-TextCell $ws4 "I42" "//i1.hdslb.com/bfs/openplatform/202403/lqXD63661711623533572.png"

# Row 43 <- old row 42 (杭州·第八届YH樱花动漫游戏文化节), F bumped 872 -> 875
Set-TextCell $ws4 "B43" "2024-06-08"
Set-TextCell $ws4 "C43" "杭州·第八届YH樱花动漫游戏文化节"
Set-TextCell $ws4 "D43" "德胜东路2539号 梦马汽车小镇"
Set-TextCell $ws4 "E43" "2024.06.08 10:00-06.10 17:00"
Set-NumCell  $ws4 "F43" 875
Set-NumCell  $ws4 "G43" 65
Set-TextCell $ws4 "H43" "https://show.bilibili.com/platform/detail.html?id=82687"
Set-TextCell $ws4 "I43" "//i2.hdslb.com/bfs/openplatform/202403/S5pnadXj1710210939138.png"

# Row 44 <- old row 43 (杭州·第三届日夜国乙only), F bumped 1582 -> 1584
# (old row 44, "杭州·苗阜王声 青曲社相声全国巡演", is dropped by the shift)
Set-TextCell $ws4 "B44" "2024-06-09"
Set-TextCell $ws4 "C44" "杭州·第三届日夜国乙only"
Set-TextCell $ws4 "D44" "创意路1号 中国智谷富春园区"
Set-TextCell $ws4 "E44" "2024.06.09 10:00-06.09 23:00"
Set-NumCell  $ws4 "F44" 1584
Set-NumCell  $ws4 "G44" 58
Set-TextCell $ws4 "H44" "https://show.bilibili.com/platform/detail.html?id=82618"
Set-TextCell $ws4 "I44" "//i2.hdslb.com/bfs/openplatform/202403/fXRzYEFH1710124366279.png"

# Row 45 keeps its event (杭州·代号鸢only-广陵大学); only the F count changes.
Set-NumCell $ws4 "F45" 266

# Trailing simple F-count bumps (unaffected by the row shift above).
Set-NumCell $ws4 "F47" 109
Set-NumCell $ws4 "F48" 97
Set-NumCell $ws4 "F49" 136
